$d = $word.ActiveDocument

# 1) Update the revision date/time stamp.
$d.Content.Find.Execute(
    "May  31, 2021 (06:10:01 PM)", $false, $false, $false, $false, $false,
    $true, 1, $false, "May  31, 2021 (06:38:15 PM)", 2)

# 2) Tweak the intro sentence: "for how to create a class in Visual Studio"
#    becomes "how to create a class in your IDE".
$d.Content.Find.Execute(
    "and for how to create a class in Visual Studio", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "and how to create a class in your IDE", 2)

# 3) Insert a brand-new bullet right after the "Comment out the body of the
#    Main method..." bullet, before "Your program should compile...".
#    We anchor the Find on the trailing closing curly-quote + period of that
#    bullet (scoped to that single paragraph) and replace it with itself
#    plus a paragraph break and the new bullet text, so the new paragraph
#    naturally inherits the same paragraph formatting (Compact style /
#    numbering) without picking up stray character styles.
$p = $d.Paragraphs(22)
$r = $p.Range
$closeQuote = [char]8221
$newBulletText = "It is important that you re-name the files within your IDE. If you try to rename your files, or their folders, outside of the IDE then it will break your solution. The solution will still be looking for the original file/folder names, and will not recognize the changed names. If such an error occurs, restore the previous names and then rename your files through the IDE as instructed."
$r.Find.Execute(
    $closeQuote + ".", $false, $false, $false, $false, $false, $true, 1,
    $false, $closeQuote + "." + [char]13 + $newBulletText, 2)
